# Contrat_CDI.docx - Article 19 clause title cleanup
#
# The paragraph that used to read
#     "Article 19 : Clause de Protection des données à caractère personnel"
# loses the leading "Article 19 : " label, leaving only
#     "Clause de Protection des données à caractère personnel"
# Deleting that text in Word moves the (hidden) "_GoBack" last-edit
# bookmark to the point of the edit, so we relocate it from its old
# spot (in the middle of the "${date_redaction}" merge field near the
# end of the document) to right at the start of the now-shortened
# paragraph.

$d = $word.ActiveDocument

# 1) Delete the "Article 19 : " run (wdFindContinue = 1, wdReplaceAll = 2).
$null = $d.Content.Find.Execute(
    "Article 19 : ",  # FindText
    $false,           # MatchCase
    $false,           # MatchWholeWord
    $false,           # MatchWildcards
    $false,           # MatchSoundsLike
    $false,           # MatchAllWordForms
    $true,            # Forward
    1,                # Wrap (wdFindContinue)
    $false,           # Format
    "",               # ReplaceWith
    2                 # Replace (wdReplaceAll)
)

# 2) Re-point "_GoBack" at the start of "Clause de Protection...".
#    Adding a bookmark named "_GoBack" behaves like Word's own
#    last-edit-position bookmark: it is unique, so it is simply moved
#    off its previous location onto the new one.
$target = $d.Content
$null = $target.Find.Execute(
    "Clause de Protection des données à caractère personnel",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)
